$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 500
$ws.Range("K12").Value = 500
$ws.Range("M12").Value = -330
$ws.Range("H17").Value = 2127.1428
$ws.Range("J17").Value = 2127.1428
$ws.Range("L17").Value = 6381.428400000001
$ws.Range("N17").Value = -6717.428400000001
$ws.Range("H33").Value = 3649.5862
$ws.Range("I33").Value = 4476.609
$ws.Range("K33").Value = 4476.609
$ws.Range("M33").Value = -4247.609
$ws.Range("H57").Value = 69999
$ws.Range("J57").Value = 69999
$ws.Range("L57").Value = 209997
$ws.Range("N57").Value = -210995
$ws.Range("H76").Value = 7244.9165
$ws.Range("I76").Value = 5215.778
$ws.Range("J76").Value = 13332.333
$ws.Range("K76").Value = 5215.778
$ws.Range("L76").Value = 13332.333
$ws.Range("M76").Value = -4900.778
$ws.Range("N76").Value = -13962.333
$ws.Range("H79").Value = 7244.9165
$ws.Range("I79").Value = 5215.778
$ws.Range("J79").Value = 13332.333
$ws.Range("K79").Value = 5215.778
$ws.Range("L79").Value = 13332.333
$ws.Range("M79").Value = -4123.778
$ws.Range("N79").Value = -15516.333
$ws.Range("H103").Value = 2649.8333
$ws.Range("I103").Value = 10000
$ws.Range("J103").Value = 1179.8
$ws.Range("K103").Value = 30000
$ws.Range("L103").Value = 3539.4
$ws.Range("M103").Value = -29414
$ws.Range("N103").Value = -4711.4
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H141").Value = 919.46155
$ws.Range("I141").Value = 919.46155
$ws.Range("K141").Value = 2758.38465
$ws.Range("M141").Value = 2421.61535

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 999
$ws.Range("I10").Value = 999
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 999
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -829
$ws.Range("N10").ClearContents()
$ws.Range("H97").Value = 5524.5
$ws.Range("J97").Value = 3125.9
$ws.Range("L97").Value = 3125.9
$ws.Range("N97").Value = -4117.9
$ws.Range("H132").Value = 26000.791
$ws.Range("I132").Value = 34995.066
$ws.Range("K132").Value = 104985.198
$ws.Range("M132").Value = -102455.198

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3572915.5
$ws.Range("I20").Value = 4546504.5
$ws.Range("J20").Value = 3089.5
$ws.Range("K20").Value = 4546504.5
$ws.Range("L20").Value = 3089.5
$ws.Range("M20").Value = -4546257.5
$ws.Range("N20").Value = -3583.5
$ws.Range("H22").Value = 126173.875
$ws.Range("I22").Value = 143984.42
$ws.Range("K22").Value = 143984.42
$ws.Range("M22").Value = -143811.42

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1801.8
$ws.Range("I22").Value = 1533
$ws.Range("J22").Value = 1869
$ws.Range("K22").Value = 1533
$ws.Range("L22").Value = 1869
$ws.Range("M22").Value = -1183
$ws.Range("N22").Value = -2569
$ws.Range("H62").Value = 6023.5
$ws.Range("I62").Value = 5948.5
$ws.Range("J62").Value = 6048.5
$ws.Range("K62").Value = 5948.5
$ws.Range("L62").Value = 6048.5
$ws.Range("M62").Value = -5324.5
$ws.Range("N62").Value = -7296.5
$ws.Range("H65").Value = 6023.5
$ws.Range("I65").Value = 5948.5
$ws.Range("J65").Value = 6048.5
$ws.Range("K65").Value = 29742.5
$ws.Range("L65").Value = 30242.5
$ws.Range("M65").Value = -26622.5
$ws.Range("N65").Value = -36482.5
$ws.Range("H86").Value = 4964.8
$ws.Range("I86").Value = 4895
$ws.Range("K86").Value = 4895
$ws.Range("M86").Value = -3772
$ws.Range("H89").Value = 4964.8
$ws.Range("I89").Value = 4895
$ws.Range("K89").Value = 24475
$ws.Range("M89").Value = -18859
$ws.Range("H103").Value = 35003.668
$ws.Range("I103").Value = 17755.75
$ws.Range("J103").Value = 69499.5
$ws.Range("K103").Value = 17755.75
$ws.Range("L103").Value = 69499.5
$ws.Range("M103").Value = -16583.75
$ws.Range("N103").Value = -71843.5
$ws.Range("H105").Value = 2126.6191
$ws.Range("I105").Value = 1957.6111
$ws.Range("J105").Value = 3140.6667
$ws.Range("K105").Value = 1957.6111
$ws.Range("L105").Value = 3140.6667
$ws.Range("M105").Value = -210.6111000000001
$ws.Range("N105").Value = -6634.6667
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H107").Value = 672.4667
$ws.Range("I107").Value = 672.4667
$ws.Range("K107").Value = 672.4667
$ws.Range("M107").Value = 1247.5333

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 427
$ws.Range("I29").Value = 262.05
$ws.Range("J29").Value = 610.2778
$ws.Range("K29").Value = 786.1500000000001
$ws.Range("L29").Value = 1830.8334
$ws.Range("M29").Value = -509.1500000000001
$ws.Range("N29").Value = -2384.8334
$ws.Range("H34").Value = 995.4167
$ws.Range("J34").Value = 1000
$ws.Range("L34").Value = 3000
$ws.Range("N34").Value = -3168
$ws.Range("H88").Value = 10023
$ws.Range("J88").Value = 14508
$ws.Range("L88").Value = 43524
$ws.Range("N88").Value = -44380
$ws.Range("H91").Value = 10023
$ws.Range("J91").Value = 14508
$ws.Range("L91").Value = 43524
$ws.Range("N91").Value = -46488
$ws.Range("H129").Value = 719510
$ws.Range("I129").Value = 16583.857
$ws.Range("J129").Value = 1703606.6
$ws.Range("K129").Value = 49751.571
$ws.Range("L129").Value = 5110819.800000001
$ws.Range("M129").Value = -44751.571
$ws.Range("N129").Value = -5120819.800000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 2985.6667
$ws.Range("I13").Value = 2985.6667
$ws.Range("K13").Value = 2985.6667
$ws.Range("M13").Value = -2846.6667
$ws.Range("H15").Value = 49999.332
$ws.Range("J15").Value = 49999.332
$ws.Range("L15").Value = 49999.332
$ws.Range("N15").Value = -50575.332
$ws.Range("H81").Value = 49999.332
$ws.Range("J81").Value = 49999.332
$ws.Range("L81").Value = 49999.332
$ws.Range("N81").Value = -51995.332
$ws.Range("H84").Value = 49999.332
$ws.Range("J84").Value = 49999.332
$ws.Range("L84").Value = 149997.996
$ws.Range("N84").Value = -159981.996

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 19997.5
$ws.Range("J12").Value = 19997.5
$ws.Range("L12").Value = 19997.5
$ws.Range("N12").Value = -20337.5
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 10000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -9826
$ws.Range("N21").ClearContents()
$ws.Range("H100").Value = 3667
$ws.Range("I100").Value = 3668.3
$ws.Range("K100").Value = 3668.3
$ws.Range("M100").Value = -3127.3
$ws.Range("H136").Value = 1841.381
$ws.Range("I136").Value = 1743.45
$ws.Range("K136").Value = 5230.35
$ws.Range("M136").Value = -2680.35

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 49997.5
$ws.Range("J86").Value = 49997.5
$ws.Range("L86").Value = 49997.5
$ws.Range("N86").Value = -52243.5
$ws.Range("H89").Value = 49997.5
$ws.Range("J89").Value = 49997.5
$ws.Range("L89").Value = 249987.5
$ws.Range("N89").Value = -261219.5
$ws.Range("H110").Value = 79999
$ws.Range("J110").Value = 79999
$ws.Range("L110").Value = 79999
$ws.Range("N110").Value = -88179
$ws.Range("H116").Value = 84840
$ws.Range("J116").Value = 84840
$ws.Range("L116").Value = 84840
$ws.Range("N116").Value = -94018
